# Fix Training Data Issue
# The "Date" column (BF) held a malformed label ("6-30-2011-12") on every
# data row. Replace it with the correct ISO-style date string "2012-06-30".
# NumberFormat is forced to Text ("@") first so Excel stores the value as a
# literal string instead of re-interpreting/recalculating it as a date serial.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateCol = 58   # column BF
$firstRow = 2
$lastRow = 31
$newValue = "2012-06-30"

$targetRange = $ws.Range($ws.Cells.Item($firstRow, $dateCol), $ws.Cells.Item($lastRow, $dateCol))
$targetRange.NumberFormat = "@"

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, $dateCol).Value = $newValue
}
